{"js": "// Apply the three changes described by the diff:\n// 1. \"last updated:\" date changes from 2022-06-28 to 2022-11-04.\n// 2. Six \"## Warning: package '...' was built under R version ...\" source-code\n//    paragraphs (right after the \"palmerpenguins\" sentence, before the\n//    \"Motivation\" section's bookmarkEnd) are removed entirely.\n// 3. \"## # A tibble: 3 x 6\" becomes \"## # A tibble: 3 \u00d7 6\" (ascii x -> \u00d7).\n\nconst body = context.document.body;\n\n// --- Change 1: update the \"last updated\" date -----------------------------\nconst dateResults = body.search(\"2022-06-28\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].insertText(\"2022-11-04\", \"Replace\");\n}\nawait context.sync();\n\n// --- Change 2: delete the six R-warning paragraphs -------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst warningParagraphs = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"was built under R version\") >= 0) {\n    warningParagraphs.push(paragraphs.items[i]);\n  }\n}\nfor (const p of warningParagraphs) {\n  p.delete();\n}\nawait context.sync();\n\n// --- Change 3: replace \"x\" with \"\u00d7\" in the tibble dimension line -----------\nconst tibbleResults = body.search(\"A tibble: 3 x 6\", { matchCase: true });\ntibbleResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < tibbleResults.items.length; i++) {\n  tibbleResults.items[i].insertText(\"A tibble: 3 \u00d7 6\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Apply the three changes described by the diff:\n# 1. \"last updated:\" date changes from 2022-06-28 to 2022-11-04.\n# 2. Six \"## Warning: package '...' was built under R version ...\" source-code\n#    paragraphs (right after the \"palmerpenguins\" sentence, before the\n#    \"Motivation\" section's bookmarkEnd) are removed entirely.\n# 3. \"## # A tibble: 3 x 6\" becomes \"## # A tibble: 3 \u00d7 6\" (ascii x -> \u00d7).\n\n$d = $word.ActiveDocument\n\n# --- Change 1: update the \"last updated\" date ------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"2022-06-28\"\n$find.Replacement.Text = \"2022-11-04\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# --- Change 2: delete the six R-warning paragraphs --------------------------\n$count = $d.Paragraphs.Count\n$toDelete = @()\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -like \"*was built under R version*\") {\n        $toDelete += $i\n    }\n}\nfor ($j = $toDelete.Count - 1; $j -ge 0; $j--) {\n    $idx = $toDelete[$j]\n    $p = $d.Paragraphs.Item($idx)\n    $p.Range.Delete()\n}\n\n# --- Change 3: replace \"x\" with \"\u00d7\" in the tibble dimension line ------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"A tibble: 3 x 6\"\n$find2.Replacement.Text = \"A tibble: 3 \u00d7 6\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
